$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry was added at the top of the data block (row 51),
# pushing every existing record down by one row (51:136 -> 52:137).
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the latest weekly entry.
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44967
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112001
$ws.Range("G51").Value = "Berenjena"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 8500
$ws.Range("M51").Value = 8250
$ws.Range("N51").Value = "$/caja 60 unidades"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 138
$ws.Range("Q51").Value = 60
$ws.Range("R51").Value = "Hortaliza"
